$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Eventos")

$ws.Range("K2").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743717012329/Ggg0CQFNgb-camera-video-segment-driver-1743717014829.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCCSGOIM3%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T170837Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCzJ%2BsNb9QYG5jEVuoyfi7%2FGDu79tW1SKiHEhHptlN8PAIgZH1ThIJsUl44fmtU2AGlH0TEEKJZ%2FTyf3dECJGRY5lwq3QMIFBAEGgw3ODEyMDQ5NDIyNDQiDESm0Aw31mMQfpQsZiq6A70pCSetwaD7Cug43hBZY4aA9Kyd5V6dd%2FBzCu9Y8FB4ykVWCmxVDmxGuYr3n6T7LXF%2FJTzZq4ewN5i%2BCGKh4FuHzzaGHHL5QHsfoaj5GD2RKAx698tBzrVyFFAQQBf7oryfYaCGOAoooPEEfJVg6CS9gKXIgMar099z2bUtSNoahSzyMF1F1vYOHLKJCbmcg99SzQHtfs7Fr44w92c6lTW26jUzIrk%2B7lpCMRklm5D7LNG%2B4PDpuBCErlFHhI8QR7OP%2FgoLS1JQswx3doad1E%2BNEOcJI7q3BJV1CnOIdyuNjeK31JSgueZgZS5Mmz%2BT%2B1MdY2fzpxAxPcXI4VLdsOtz237mpx8RokMjKPVkxAXN2dwWU5Tdt6LWYyWyPU49TV8c%2FJm1D8t%2FLjHy5wfXhSUd8Jd17VxxP3iGgu36uz9pIcW0AgyP8r0e5gSYbHttCYBgIPwIMFAD4bNRsNCAOjtTTd3AKDCQPCI%2F6tSsjhsMQrG8vGDxv3BHNtaefCJ54Wco01knhQK88M8q7XKjBFZI%2FdaKfV6fd59s6FYqOncalcTHKeJ2WAV9c5JpLsjMdRcKIyDCBjH8ax8w5P2%2BvwY6pQFqWMuhsFheb1rD7dABPXyBAUC7Iy7wgOwrLEfySv74NBtJGX2MImfHR3shtNcMosVXBmPkOqTQzIgNTWcBMrJw9umyK9B0g6oSPU4EwuhlW7NQCLEwyTMFJMzq2n5Z1%2FLHnfsIustaCH%2FoUMbw7btHbQL0%2BWdil%2BCXDg2QgqBx265ejmg4RbYxfN3M7XO2s6ofpnm1UoZTxiJn%2FKu5WrZOl4CNR2Y%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2001%3A08%3A37%20GMT&X-Amz-Signature=e96325253753611fcff5d21f4de181439661f10d99c2bb72941ed40d9d20ad29'
$ws.Range("K3").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743713574485/FgwSeOefdQ-camera-video-segment-driver-1743713576985.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCCSGOIM3%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T170837Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCzJ%2BsNb9QYG5jEVuoyfi7%2FGDu79tW1SKiHEhHptlN8PAIgZH1ThIJsUl44fmtU2AGlH0TEEKJZ%2FTyf3dECJGRY5lwq3QMIFBAEGgw3ODEyMDQ5NDIyNDQiDESm0Aw31mMQfpQsZiq6A70pCSetwaD7Cug43hBZY4aA9Kyd5V6dd%2FBzCu9Y8FB4ykVWCmxVDmxGuYr3n6T7LXF%2FJTzZq4ewN5i%2BCGKh4FuHzzaGHHL5QHsfoaj5GD2RKAx698tBzrVyFFAQQBf7oryfYaCGOAoooPEEfJVg6CS9gKXIgMar099z2bUtSNoahSzyMF1F1vYOHLKJCbmcg99SzQHtfs7Fr44w92c6lTW26jUzIrk%2B7lpCMRklm5D7LNG%2B4PDpuBCErlFHhI8QR7OP%2FgoLS1JQswx3doad1E%2BNEOcJI7q3BJV1CnOIdyuNjeK31JSgueZgZS5Mmz%2BT%2B1MdY2fzpxAxPcXI4VLdsOtz237mpx8RokMjKPVkxAXN2dwWU5Tdt6LWYyWyPU49TV8c%2FJm1D8t%2FLjHy5wfXhSUd8Jd17VxxP3iGgu36uz9pIcW0AgyP8r0e5gSYbHttCYBgIPwIMFAD4bNRsNCAOjtTTd3AKDCQPCI%2F6tSsjhsMQrG8vGDxv3BHNtaefCJ54Wco01knhQK88M8q7XKjBFZI%2FdaKfV6fd59s6FYqOncalcTHKeJ2WAV9c5JpLsjMdRcKIyDCBjH8ax8w5P2%2BvwY6pQFqWMuhsFheb1rD7dABPXyBAUC7Iy7wgOwrLEfySv74NBtJGX2MImfHR3shtNcMosVXBmPkOqTQzIgNTWcBMrJw9umyK9B0g6oSPU4EwuhlW7NQCLEwyTMFJMzq2n5Z1%2FLHnfsIustaCH%2FoUMbw7btHbQL0%2BWdil%2BCXDg2QgqBx265ejmg4RbYxfN3M7XO2s6ofpnm1UoZTxiJn%2FKu5WrZOl4CNR2Y%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2001%3A08%3A37%20GMT&X-Amz-Signature=cf2db0dedde93f48200ca0229d5bf679fea8e8abba66f9784f31577c17bdda31'
$ws.Range("K4").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712999664/UzzO0YQTJX-camera-video-segment-driver-1743713002164.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCCSGOIM3%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T170837Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCzJ%2BsNb9QYG5jEVuoyfi7%2FGDu79tW1SKiHEhHptlN8PAIgZH1ThIJsUl44fmtU2AGlH0TEEKJZ%2FTyf3dECJGRY5lwq3QMIFBAEGgw3ODEyMDQ5NDIyNDQiDESm0Aw31mMQfpQsZiq6A70pCSetwaD7Cug43hBZY4aA9Kyd5V6dd%2FBzCu9Y8FB4ykVWCmxVDmxGuYr3n6T7LXF%2FJTzZq4ewN5i%2BCGKh4FuHzzaGHHL5QHsfoaj5GD2RKAx698tBzrVyFFAQQBf7oryfYaCGOAoooPEEfJVg6CS9gKXIgMar099z2bUtSNoahSzyMF1F1vYOHLKJCbmcg99SzQHtfs7Fr44w92c6lTW26jUzIrk%2B7lpCMRklm5D7LNG%2B4PDpuBCErlFHhI8QR7OP%2FgoLS1JQswx3doad1E%2BNEOcJI7q3BJV1CnOIdyuNjeK31JSgueZgZS5Mmz%2BT%2B1MdY2fzpxAxPcXI4VLdsOtz237mpx8RokMjKPVkxAXN2dwWU5Tdt6LWYyWyPU49TV8c%2FJm1D8t%2FLjHy5wfXhSUd8Jd17VxxP3iGgu36uz9pIcW0AgyP8r0e5gSYbHttCYBgIPwIMFAD4bNRsNCAOjtTTd3AKDCQPCI%2F6tSsjhsMQrG8vGDxv3BHNtaefCJ54Wco01knhQK88M8q7XKjBFZI%2FdaKfV6fd59s6FYqOncalcTHKeJ2WAV9c5JpLsjMdRcKIyDCBjH8ax8w5P2%2BvwY6pQFqWMuhsFheb1rD7dABPXyBAUC7Iy7wgOwrLEfySv74NBtJGX2MImfHR3shtNcMosVXBmPkOqTQzIgNTWcBMrJw9umyK9B0g6oSPU4EwuhlW7NQCLEwyTMFJMzq2n5Z1%2FLHnfsIustaCH%2FoUMbw7btHbQL0%2BWdil%2BCXDg2QgqBx265ejmg4RbYxfN3M7XO2s6ofpnm1UoZTxiJn%2FKu5WrZOl4CNR2Y%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2001%3A08%3A37%20GMT&X-Amz-Signature=e7dfa7950f8db12c8032a3c7b0e368470de42dbf91a0753be97f659b422b493a'
$ws.Range("K5").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991109864/1743712546260/Jh2p7rn04n-camera-video-segment-driver-1743712548760.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCCSGOIM3%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T170837Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCzJ%2BsNb9QYG5jEVuoyfi7%2FGDu79tW1SKiHEhHptlN8PAIgZH1ThIJsUl44fmtU2AGlH0TEEKJZ%2FTyf3dECJGRY5lwq3QMIFBAEGgw3ODEyMDQ5NDIyNDQiDESm0Aw31mMQfpQsZiq6A70pCSetwaD7Cug43hBZY4aA9Kyd5V6dd%2FBzCu9Y8FB4ykVWCmxVDmxGuYr3n6T7LXF%2FJTzZq4ewN5i%2BCGKh4FuHzzaGHHL5QHsfoaj5GD2RKAx698tBzrVyFFAQQBf7oryfYaCGOAoooPEEfJVg6CS9gKXIgMar099z2bUtSNoahSzyMF1F1vYOHLKJCbmcg99SzQHtfs7Fr44w92c6lTW26jUzIrk%2B7lpCMRklm5D7LNG%2B4PDpuBCErlFHhI8QR7OP%2FgoLS1JQswx3doad1E%2BNEOcJI7q3BJV1CnOIdyuNjeK31JSgueZgZS5Mmz%2BT%2B1MdY2fzpxAxPcXI4VLdsOtz237mpx8RokMjKPVkxAXN2dwWU5Tdt6LWYyWyPU49TV8c%2FJm1D8t%2FLjHy5wfXhSUd8Jd17VxxP3iGgu36uz9pIcW0AgyP8r0e5gSYbHttCYBgIPwIMFAD4bNRsNCAOjtTTd3AKDCQPCI%2F6tSsjhsMQrG8vGDxv3BHNtaefCJ54Wco01knhQK88M8q7XKjBFZI%2FdaKfV6fd59s6FYqOncalcTHKeJ2WAV9c5JpLsjMdRcKIyDCBjH8ax8w5P2%2BvwY6pQFqWMuhsFheb1rD7dABPXyBAUC7Iy7wgOwrLEfySv74NBtJGX2MImfHR3shtNcMosVXBmPkOqTQzIgNTWcBMrJw9umyK9B0g6oSPU4EwuhlW7NQCLEwyTMFJMzq2n5Z1%2FLHnfsIustaCH%2FoUMbw7btHbQL0%2BWdil%2BCXDg2QgqBx265ejmg4RbYxfN3M7XO2s6ofpnm1UoZTxiJn%2FKu5WrZOl4CNR2Y%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2001%3A08%3A37%20GMT&X-Amz-Signature=33b88a7d34f7770eac0ffae3f0e064f822c00127fde92a98e2afff3e7a751a91'
$ws.Range("K6").Value = 'https://s3.samsara.com/samsara-cvdata/4006124/281474991205262/1743703000742/spTbAbfADq-camera-video-segment-driver-1743703003242.audio.mp4?X-Amz-Algorithm=AWS4-HMAC-SHA256&X-Amz-Credential=ASIA3LY3RNWSCCSGOIM3%2F20250404%2Fus-west-2%2Fs3%2Faws4_request&X-Amz-Date=20250404T170837Z&X-Amz-Expires=28800&X-Amz-Security-Token=IQoJb3JpZ2luX2VjEJv%2F%2F%2F%2F%2F%2F%2F%2F%2F%2FwEaCXVzLXdlc3QtMiJHMEUCIQCzJ%2BsNb9QYG5jEVuoyfi7%2FGDu79tW1SKiHEhHptlN8PAIgZH1ThIJsUl44fmtU2AGlH0TEEKJZ%2FTyf3dECJGRY5lwq3QMIFBAEGgw3ODEyMDQ5NDIyNDQiDESm0Aw31mMQfpQsZiq6A70pCSetwaD7Cug43hBZY4aA9Kyd5V6dd%2FBzCu9Y8FB4ykVWCmxVDmxGuYr3n6T7LXF%2FJTzZq4ewN5i%2BCGKh4FuHzzaGHHL5QHsfoaj5GD2RKAx698tBzrVyFFAQQBf7oryfYaCGOAoooPEEfJVg6CS9gKXIgMar099z2bUtSNoahSzyMF1F1vYOHLKJCbmcg99SzQHtfs7Fr44w92c6lTW26jUzIrk%2B7lpCMRklm5D7LNG%2B4PDpuBCErlFHhI8QR7OP%2FgoLS1JQswx3doad1E%2BNEOcJI7q3BJV1CnOIdyuNjeK31JSgueZgZS5Mmz%2BT%2B1MdY2fzpxAxPcXI4VLdsOtz237mpx8RokMjKPVkxAXN2dwWU5Tdt6LWYyWyPU49TV8c%2FJm1D8t%2FLjHy5wfXhSUd8Jd17VxxP3iGgu36uz9pIcW0AgyP8r0e5gSYbHttCYBgIPwIMFAD4bNRsNCAOjtTTd3AKDCQPCI%2F6tSsjhsMQrG8vGDxv3BHNtaefCJ54Wco01knhQK88M8q7XKjBFZI%2FdaKfV6fd59s6FYqOncalcTHKeJ2WAV9c5JpLsjMdRcKIyDCBjH8ax8w5P2%2BvwY6pQFqWMuhsFheb1rD7dABPXyBAUC7Iy7wgOwrLEfySv74NBtJGX2MImfHR3shtNcMosVXBmPkOqTQzIgNTWcBMrJw9umyK9B0g6oSPU4EwuhlW7NQCLEwyTMFJMzq2n5Z1%2FLHnfsIustaCH%2FoUMbw7btHbQL0%2BWdil%2BCXDg2QgqBx265ejmg4RbYxfN3M7XO2s6ofpnm1UoZTxiJn%2FKu5WrZOl4CNR2Y%3D&X-Amz-SignedHeaders=host&response-expires=Sat%2C%2005%20Apr%202025%2001%3A08%3A37%20GMT&X-Amz-Signature=6204e043d1a0bebbf30f5b01fb7a7fa2e3ca5ade5e31b9c2250a4f39091030b7'
